$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.671.19'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '2.277.85'
$ws.Range("E3").Value = '  -0.60%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.38%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '113.56'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +9.80%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '266.96'
$c.NumberFormat = "General"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.624'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.88%  '
$ws.Range("E8").Value = '  +0.31%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.609'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.10%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '48.22'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +4.82%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0934'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.20%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '8.81'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +7.81%  '
$ws.Range("E13").Value = '  +0.34%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '15.67'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.79%  '
$ws.Range("D15").Value = '2.624.08'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.868'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").Value = '2.278.42'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").Value = '43.489.73'
$ws.Range("E18").Value = '  -0.62%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.0000108'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.53%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.02'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +11.71%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '71.87'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.59%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.40'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -4.47%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.95'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +7.36%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '232.08'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("E26").Value = '  -0.01%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '11.47'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.45%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '41.07'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.04%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '3.39'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.59%  '
$ws.Range("E30").Value = '  +1.25%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '173.14'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.62%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '21.42'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("E33").Value = '  +0.77%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '5.63'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.80%  '
$ws.Range("E35").Value = '  -0.04%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '4.59'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = '  -6.20%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.0350'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.41%  '
$ws.Range("E38").Value = '  -5.30%  '
$ws.Range("E39").Value = '  +5.54%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '14.34'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = '  +17.63%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '74.69'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +13.92%  '
$ws.Range("E42").Value = '  +3.44%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.236'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.19%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '6.19'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +16.72%  '
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("E46").Value = '  +0.10%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '8.65'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.77%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0995'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.25'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '101.62'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +2.36%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.452'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.96%  '
